$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.910.78'
$ws.Range('E2').Value = '  -2.09%  '
$ws.Range('D3').Value = '1.899.88'
$ws.Range('E3').Value = '  -4.11%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'324.18"
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('D7').Value = "'0.4581"
$ws.Range('E7').Value = '  -1.87%  '
$ws.Range('E8').Value = '  -2.68%  '
$ws.Range('D9').Value = "'0.07709"
$ws.Range('E9').Value = '  -3.26%  '
$ws.Range('D10').Value = "'0.9747"
$ws.Range('E10').Value = '  -2.04%  '
$ws.Range('D11').Value = "'21.98"
$ws.Range('E11').Value = '  -4.04%  '
$ws.Range('D12').Value = '1.897.03'
$ws.Range('E12').Value = '  -4.81%  '
$ws.Range('D13').Value = "'6.919"
$ws.Range('E13').Value = '  -4.10%  '
$ws.Range('D14').Value = "'5.633"
$ws.Range('E14').Value = '  -3.69%  '
$ws.Range('D15').Value = "'0.07030"
$ws.Range('E15').Value = '  -0.74%  '
$ws.Range('D16').Value = "'1.004"
$ws.Range('E16').Value = '  -0.19%  '
$ws.Range('D17').Value = "'83.55"
$ws.Range('E17').Value = '  -4.97%  '
$ws.Range('D18').Value = "'0.000009468"
$ws.Range('E18').Value = '  -4.95%  '
$ws.Range('D19').Value = "'16.59"
$ws.Range('E19').Value = '  -4.30%  '
$ws.Range('D20').Value = "'1.002"
$ws.Range('E20').Value = '  -0.29%  '
$ws.Range('D21').Value = '28.889.49'
$ws.Range('E21').Value = '  -2.12%  '
$ws.Range('D22').Value = "'5.280"
$ws.Range('E22').Value = '  -5.27%  '
$ws.Range('E23').Value = '  -3.35%  '
$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D24').Value = '2.119.98'
$ws.Range('E24').Value = '  -5.13%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = "'2.094"
$ws.Range('E25').Value = '  -0.74%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = "'158.12"
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = "'18.99"
$ws.Range('E27').Value = '  -3.21%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = "'5.603"
$ws.Range('E28').Value = '  -3.50%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = "'117.28"
$ws.Range('E29').Value = '  -2.09%  '
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D30').Value = "'1.833"
$ws.Range('E30').Value = '  -4.12%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = "'0.09231"
$ws.Range('E31').Value = '  -2.10%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = "'0.8588"
$ws.Range('E32').Value = '  -4.50%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'5.072"
$ws.Range('E33').Value = '  -3.29%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = "'1.236"
$ws.Range('E34').Value = '  -6.94%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = "'2.987"
$ws.Range('E35').Value = '  -6.48%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = "'0.05654"
$ws.Range('E36').Value = '  -2.90%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = "'1.139"
$ws.Range('E37').Value = '  -2.68%  '
$ws.Range('B38').Value = 'Frax'
$ws.Range('C38').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D38').Value = "'1.002"
$ws.Range('E38').Value = '  -0.13%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = "'0.02032"
$ws.Range('E39').Value = '  -3.38%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = "'0.5468"
$ws.Range('E40').Value = '  -4.64%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = "'7.371"
$ws.Range('E41').Value = '  -5.53%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = "'0.1747"
$ws.Range('E42').Value = '  -3.08%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = "'9.248"
$ws.Range('E43').Value = '  -4.39%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = "'2.750"
$ws.Range('E44').Value = '  -1.44%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = "'0.5138"
$ws.Range('E45').Value = '  -4.45%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = "'11.16"
$ws.Range('E46').Value = '  -5.57%  '
$ws.Range('D47').Value = "'0.000002628"
$ws.Range('E47').Value = '  -15.93%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = "'0.06806"
$ws.Range('E48').Value = '  -1.77%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = "'2.066"
$ws.Range('E49').Value = '  -5.63%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = "'109.85"
$ws.Range('E50').Value = '  -3.89%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = "'1.765"
$ws.Range('E51').Value = '  -3.68%  '
